# "Added Darley Abbey run"
# Updates the "start/end at pub" type to "start/end at pub + extra pub" for
# the multi-pub runs, marks the Duffield run's PM column as "nr", and adds
# the new Darley Abbey (The Abbey / The Furnace) run as row 79.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NB: new shared strings are interned in first-seen order, so touch the
# cells in the same order the author's saved file introduces them:
#   "The Abbey / The Furnace" -> "start/end at pub + extra pub"
#   -> "The Abbey basement" -> "Drunkard taxi drama..."

# 3) New run: Darley Abbey, 2024-10-23 -- The Abbey / The Furnace.
$ws.Range("A79").Value = 45588
$ws.Range("B79").Value = "The Abbey / The Furnace"
$ws.Range("C79").Value = "Darley Abbey"

# 1) Existing multi-pub-run rows get the new, more specific "Type" label.
$ws.Range("D34").Value = "start/end at pub + extra pub"
$ws.Range("D51").Value = "start/end at pub + extra pub"
$ws.Range("D52").Value = "start/end at pub + extra pub"
$ws.Range("D62").Value = "start/end at pub + extra pub"
$ws.Range("D78").Value = "start/end at pub + extra pub"
$ws.Range("D79").Value = "start/end at pub + extra pub"

$ws.Range("E79").Value = 3.57
$ws.Range("F79").Formula = "=E79*0.6213712"
$ws.Range("G79").Value = 0.02011574074074074
$ws.Range("H79").Formula = "=G79/F79"
$ws.Range("I79").Value = 1

# 2) Row 78 (Duffield, 2024-08-24) PM column was a plain "1" -- recode it to
#    "nr" like the new row below, which drops it out of the Q-column SUM.
$ws.Range("J78").Value = "nr"
$ws.Range("J79").Value = "nr"

$ws.Range("K79").Value = 1
$ws.Range("M79").Value = 1
$ws.Range("O79").Value = 1
$ws.Range("R79").Value = "The Abbey basement"
$ws.Range("P79").Value = "Drunkard taxi drama, pint at The Furnace, talkative chap at the abbey, the summoning of Mr. McCoy!"
$ws.Range("Q79").Formula = "=SUM(I79:O79)*F79"

# 4) Totals row 83 -- extend the mileage sum to cover the new row.
$ws.Range("F83").Formula = "=SUM(F8:F79)"

# 5) View state -- park the selection/scroll roughly where the author left it.
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 62
$ws.Range("P80").Select()
